# LOGBOEK.xlsx update
# Adds four new logboek entries (rows 16-19, column C) describing work done,
# matching the accompanying "Update logboek & Technical and functional
# analysis - Accelerometer Thread - Matter V02.docx" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Row 16 (day 14) ---------------------------------------------------
$ws.Range("C16").Value = "Testen dat ik heb uitgevoerd gedocumenteerd. Verdere documentatie over het Matter protocol"
$ws.Range("C16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 17.25

# --- Row 17 (day 15) ---------------------------------------------------
$ws.Range("C17").Value = "De testen level-control en accelerometer willen combineren in één project. Tussentijds een update van de SDK gedaan. Bij het debuggen kreeg ik de volgende error: WARNING: Failed to read memory @ address 0xFFFFFFFE. Ik kreeg dit niet opgelost."
$ws.Range("C17").WrapText = $true
$ws.Range("C17").VerticalAlignment = -4160
$ws.Rows.Item(17).RowHeight = 39

# --- Row 18 (day 16) ---------------------------------------------------
$ws.Range("C18").Value = "IDE opnieuw geinstalleerd om probleem te verhelpen maar is niet gelukt."
$ws.Range("C18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 24.75

# --- Row 19 (day 17) ---------------------------------------------------
$ws.Range("C19").Value = "Verder gewerkt aan technische analyse sectie Matter Interaction model. Probleem van de mcu is opgelost, er was een probleem met de bootloader. De oplossing was door een project te gebruiken met een external bootloader, eerst de bootloader geflashed en daarna de firmware. Nu werkt het terug."
$ws.Range("C19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 45

# --- View / selection state --------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("G19").Select()
